# Auto-generated edit script applying the Alexander_Profits.xlsx data refresh diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 62501584
$ws.Range("I100").Value = 1611.6666
$ws.Range("J100").Value = 250001500
$ws.Range("K100").Value = 1611.6666
$ws.Range("L100").Value = 250001500
$ws.Range("M100").Value = -1070.6666
$ws.Range("N100").Value = -250002582
$ws.Range("H107").Value = 1048.1305
$ws.Range("I107").Value = 1310.4117
$ws.Range("J107").Value = 305
$ws.Range("K107").Value = 1310.4117
$ws.Range("L107").Value = 305
$ws.Range("M107").Value = 609.5882999999999
$ws.Range("N107").Value = -4145
$ws.Range("H113").Value = 4331.0713
$ws.Range("I113").Value = 2480
$ws.Range("J113").Value = 5359.4443
$ws.Range("K113").Value = 2480
$ws.Range("L113").Value = 5359.4443
$ws.Range("M113").Value = 774
$ws.Range("N113").Value = -11867.4443
$ws.Range("H126").Value = 44590
$ws.Range("J126").Value = 44590
$ws.Range("L126").Value = 44590
$ws.Range("N126").Value = -54470
$ws.Range("H128").Value = 45326.668
$ws.Range("J128").Value = 45326.668
$ws.Range("L128").Value = 45326.668
$ws.Range("N128").Value = -55286.668
$ws.Range("H129").Value = 839.8125
$ws.Range("I129").Value = 521
$ws.Range("J129").Value = 1087.7778
$ws.Range("K129").Value = 1563
$ws.Range("L129").Value = 3263.3334
$ws.Range("M129").Value = 3437
$ws.Range("N129").Value = -13263.3334
$ws.Range("H133").Value = 54351.4
$ws.Range("J133").Value = 54351.4
$ws.Range("L133").Value = 54351.4
$ws.Range("N133").Value = -64471.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 23925.334
$ws.Range("J109").Value = 23925.334
$ws.Range("L109").Value = 23925.334
$ws.Range("N109").Value = -26699.334
$ws.Range("H110").Value = 5503.9287
$ws.Range("I110").Value = 5835
$ws.Range("J110").Value = 1200
$ws.Range("K110").Value = 5835
$ws.Range("L110").Value = 1200
$ws.Range("M110").Value = -3790
$ws.Range("N110").Value = -5290
$ws.Range("H122").Value = 1477.7646
$ws.Range("I122").Value = 1702.75
$ws.Range("J122").Value = 1277.7778
$ws.Range("K122").Value = 5108.25
$ws.Range("L122").Value = 3833.3334
$ws.Range("M122").Value = -2658.25
$ws.Range("N122").Value = -8733.3334
$ws.Range("H123").Value = 667375
$ws.Range("J123").Value = 667375
$ws.Range("L123").Value = 667375
$ws.Range("N123").Value = -677175
$ws.Range("H125").Value = 150044900
$ws.Range("J125").Value = 150044900
$ws.Range("L125").Value = 150044900
$ws.Range("N125").Value = -150054740
$ws.Range("H132").Value = 2124.3215
$ws.Range("I132").Value = 2023
$ws.Range("J132").Value = 2428.2856
$ws.Range("K132").Value = 6069
$ws.Range("L132").Value = 7284.8568
$ws.Range("M132").Value = -3539
$ws.Range("N132").Value = -12344.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H134").Value = 1286.2609
$ws.Range("I134").Value = 1229.2
$ws.Range("J134").Value = 1666.6666
$ws.Range("K134").Value = 3687.6
$ws.Range("L134").Value = 4999.9998
$ws.Range("M134").Value = -1152.6
$ws.Range("N134").Value = -10069.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 48996.332
$ws.Range("J20").Value = 48996.332
$ws.Range("L20").Value = 48996.332
$ws.Range("N20").Value = -49468.332
$ws.Range("H30").Value = 48996.332
$ws.Range("J30").Value = 48996.332
$ws.Range("L30").Value = 48996.332
$ws.Range("N30").Value = -49178.332
$ws.Range("H41").Value = 17737.5
$ws.Range("J41").Value = 21185
$ws.Range("L41").Value = 21185
$ws.Range("N41").Value = -22041
$ws.Range("H58").Value = 1263.4043
$ws.Range("I58").Value = 1279.0968
$ws.Range("J58").Value = 1233
$ws.Range("K58").Value = 1279.0968
$ws.Range("L58").Value = 1233
$ws.Range("M58").Value = -1076.0968
$ws.Range("N58").Value = -1639
$ws.Range("H99").Value = 1685.2273
$ws.Range("I99").Value = 1368.2858
$ws.Range("J99").Value = 2239.875
$ws.Range("K99").Value = 1368.2858
$ws.Range("L99").Value = 2239.875
$ws.Range("M99").Value = 129.7141999999999
$ws.Range("N99").Value = -5235.875
$ws.Range("H109").Value = 10928.571
$ws.Range("J109").Value = 10928.571
$ws.Range("L109").Value = 10928.571
$ws.Range("N109").Value = -13008.571
$ws.Range("H126").Value = 1685.2273
$ws.Range("I126").Value = 1368.2858
$ws.Range("J126").Value = 2239.875
$ws.Range("K126").Value = 4104.857400000001
$ws.Range("L126").Value = 6719.625
$ws.Range("M126").Value = -1634.857400000001
$ws.Range("N126").Value = -11659.625
$ws.Range("H128").Value = 48996.332
$ws.Range("J128").Value = 48996.332
$ws.Range("L128").Value = 48996.332
$ws.Range("N128").Value = -58956.332
$ws.Range("H134").Value = 4692.6665
$ws.Range("I134").Value = 5424.5835
$ws.Range("J134").Value = 1765
$ws.Range("K134").Value = 16273.7505
$ws.Range("L134").Value = 5295
$ws.Range("M134").Value = -13738.7505
$ws.Range("N134").Value = -10365
$ws.Range("H136").Value = 1263.4043
$ws.Range("I136").Value = 1279.0968
$ws.Range("J136").Value = 1233
$ws.Range("K136").Value = 3837.2904
$ws.Range("L136").Value = 3699
$ws.Range("M136").Value = -1287.2904
$ws.Range("N136").Value = -8799

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 81
$ws.Range("I38").Value = 106
$ws.Range("J38").Value = 65.375
$ws.Range("K38").Value = 318
$ws.Range("L38").Value = 196.125
$ws.Range("M38").Value = 29
$ws.Range("N38").Value = -890.125
$ws.Range("H114").Value = 414.11765
$ws.Range("I114").Value = 304.85715
$ws.Range("J114").Value = 490.6
$ws.Range("K114").Value = 914.5714499999999
$ws.Range("L114").Value = 1471.8
$ws.Range("M114").Value = 2339.42855
$ws.Range("N114").Value = -7979.8
$ws.Range("H133").Value = 2393.889
$ws.Range("I133").Value = 2130.625
$ws.Range("J133").Value = 4500
$ws.Range("K133").Value = 6391.875
$ws.Range("L133").Value = 13500
$ws.Range("M133").Value = -1331.875
$ws.Range("N133").Value = -23620
$ws.Range("H134").Value = 2041.6552
$ws.Range("I134").Value = 1439.4783
$ws.Range("J134").Value = 4350
$ws.Range("K134").Value = 4318.4349
$ws.Range("L134").Value = 13050
$ws.Range("M134").Value = 751.5650999999998
$ws.Range("N134").Value = -23190
$ws.Range("H137").Value = 2000.3478
$ws.Range("I137").Value = 1687.8572
$ws.Range("J137").Value = 2486.4443
$ws.Range("K137").Value = 5063.571599999999
$ws.Range("L137").Value = 7459.3329
$ws.Range("M137").Value = 36.42840000000069
$ws.Range("N137").Value = -17659.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 40000
$ws.Range("J127").Value = 40000
$ws.Range("L127").Value = 40000
$ws.Range("N127").Value = -49920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 54281.555
$ws.Range("J127").Value = 54281.555
$ws.Range("L127").Value = 54281.555
$ws.Range("N127").Value = -64201.555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 31510
$ws.Range("J16").Value = 31510
$ws.Range("L16").Value = 31510
$ws.Range("N16").Value = -32094
$ws.Range("H109").Value = 29438.5
$ws.Range("J109").Value = 29438.5
$ws.Range("L109").Value = 29438.5
$ws.Range("N109").Value = -32212.5
$ws.Range("H128").Value = 47245.832
$ws.Range("J128").Value = 47245.832
$ws.Range("L128").Value = 47245.832
$ws.Range("N128").Value = -57205.832
